# Apply updated dSF (column F) values to Sheet1, per "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row -> new dSF (column F) value
$updates = @{
    2  = -3
    6  = -1
    10 = -1
    11 = -3
    12 = -2
    16 = 0
    17 = 1
    21 = -1
    23 = -3
    25 = 4
    28 = -8
    29 = -1
    32 = -2
    40 = 1
    41 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
